# "control de coma por punto decimal": the numbers are entered as TEXT
# (quote-prefixed) so Excel keeps the period as the decimal mark instead of
# re-interpreting the value through the current locale's comma/point
# settings. This mirrors the source data, which grew from a 1x2 sheet
# (A1:B1) into a 4x4 table (A1:D4) of spline coefficients stored as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @("-0.0290697674418602", "0.0872093023255806", "1.02906976744186",  "0.91279069767442"),
    @("0.406976744186044",   "-3.83720930232556",  "12.8023255813953",  "-10.860465116279"),
    @("-0.686046511627902",  "9.27906976744179",   "-39.6627906976741", "59.0930232558135"),
    @("0.337209302325578",   "-6.0697674418604",   "37.0813953488368",  "-68.8139534883714")
)

for ($r = 0; $r -lt 4; $r++) {
    for ($c = 0; $c -lt 4; $c++) {
        # Leading apostrophe forces text storage (quote-prefixed), exactly
        # like a user typing '-0.029... into the cell so it is NOT parsed
        # as a locale-dependent number.
        $ws.Cells.Item($r + 1, $c + 1).Value = "'" + $values[$r][$c]
    }
}

# New header cells (C1:D1) should look like the existing A1:B1 header cells
# (bold, centered, boxed) rather than the plain default style the two new
# columns would otherwise pick up.
$ws.Range("A1:B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
